$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data cleaning: the "Authors" column (E) for rows 2-10 had author lists whose
# comma-separated entries were joined with a run of spaces. Clean the data by
# inserting one extra space into each whitespace run that follows a comma.
$rows = 2..10
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # column E = Authors
    $old = $cell.Text
    if ($old -match ',\s+') {
        $new = $old -replace ',( +)', ',$1 '
        $cell.Value = $new
    }
}
